# docs: Correctly update roadmap Excel respecting column structure
#
# Appends two new roadmap rows (7 and 8) to Sheet1, mirroring the existing
# "extra columns" (O:AB) roadmap-item block used by rows 5-6, while keeping
# the original columns (A:N) structure intact (A holds the sequential id,
# B:N stay blank placeholders).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 7: "Fourier UI Control & Window" (id 6) ----------------------
$ws.Cells.Item(7, 1).Value = 6          # A7

$ws.Cells.Item(7, 15).Value = 1                                     # O7
$ws.Cells.Item(7, 16).Value = "Fourier UI Control & Window"          # P7
$ws.Cells.Item(7, 17).Value = "UI"                                   # Q7
$ws.Cells.Item(7, 18).Value = "Aggiunto pulsante toggle per nascondere/mostrare Fourier e input numerico per definire la finestra di analisi (default 504 giorni)."  # R7
$ws.Cells.Item(7, 19).Value = "Modificato index.html (input/button), app.js (logica toggle + invio param), main.py (ricezione param), logic.py (uso param in FourierEngine)."  # S7
$ws.Cells.Item(7, 20).Value = "frontend/app.js, frontend/index.html, backend/main.py, backend/logic.py"  # T7
$ws.Cells.Item(7, 21).Value = "DONE"                                 # U7
$ws.Cells.Item(7, 22).Value = "Bassa"                                # V7
$ws.Cells.Item(7, 23).Value = "NO"                                   # W7
$ws.Cells.Item(7, 24).Value = "NO"                                   # X7
$ws.Cells.Item(7, 25).Value = "SI"                                   # Y7
$ws.Cells.Item(7, 26).Value = "Manuale"                              # Z7
# Leading apostrophe forces these to stay plain text (matching the source
# workbook, where the date columns hold literal strings, not date serials).
$ws.Cells.Item(7, 27).Value = "'2026-01-30"                          # AA7
$ws.Cells.Item(7, 28).Value = "'2026-01-30"                          # AB7

# ---- Row 8: "Media Portafoglio Equipesata" (id 7) ----------------------
$ws.Cells.Item(8, 1).Value = 7          # A8

$ws.Cells.Item(8, 15).Value = 1                                     # O8
$ws.Cells.Item(8, 16).Value = "Media Portafoglio Equipesata"         # P8
$ws.Cells.Item(8, 17).Value = "UI"                                   # Q8
$ws.Cells.Item(8, 18).Value = "Modifica calcolo media rendimenti portafoglio: ora è una media aritmetica non pesata (Equal Weighted) dei %. Ignora null."  # R8
$ws.Cells.Item(8, 19).Value = 'Refactor calcolo media in app.js. Visualizzazione "Media Equipesata".'  # S8
$ws.Cells.Item(8, 20).Value = "frontend/app.js"                      # T8
$ws.Cells.Item(8, 21).Value = "DONE"                                 # U8
$ws.Cells.Item(8, 22).Value = "Bassa"                                # V8
$ws.Cells.Item(8, 23).Value = "NO"                                   # W8
$ws.Cells.Item(8, 24).Value = "NO"                                   # X8
$ws.Cells.Item(8, 25).Value = "NO"                                   # Y8
$ws.Cells.Item(8, 26).Value = "Manuale"                              # Z8
$ws.Cells.Item(8, 27).Value = "'2026-01-30"                          # AA8
$ws.Cells.Item(8, 28).Value = "'2026-01-30"                          # AB8

# Columns B:N stay empty (inline-string placeholders) for both new rows,
# matching the pattern already used by rows 5 and 6. A bare "'" forces an
# empty TEXT cell instead of clearing the cell outright (plain "" leaves a
# blank/Number-typed cell rather than an empty string).
for ($col = 2; $col -le 14; $col++) {
    $ws.Cells.Item(7, $col).Value = "'"
    $ws.Cells.Item(8, $col).Value = "'"
}
